$d = $word.ActiveDocument

# The "capas-listado" table lists map layers; this change marks two
# finished layer rows ("Suelos" and "Carreteras") as struck-through
# (done), matching how the "Municipios" row was already formatted.
$targetLabels = @("Suelos", "Carreteras")

$t = $d.Tables.Item(1)
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $row = $t.Rows.Item($i)
    $rawLabel = $row.Cells.Item(1).Range.Text
    # Cell text carries trailing cell-mark control chars (CR + BEL);
    # strip any control characters before comparing.
    $label = $rawLabel -replace "[\x00-\x1F]", ""

    foreach ($target in $targetLabels) {
        if ($label -eq $target) {
            for ($j = 1; $j -le $row.Cells.Count; $j++) {
                $row.Cells.Item($j).Range.Font.StrikeThrough = 1
            }
        }
    }
}
